# Generate Report for Handoff
# Replaces the two "644510b7.../98644877..." files that were fully
# handed-back with three new items that are "Ready for handoff":
#   - 19a97789-...-.md   (a markdown doc, with its own handoff xlf)
#   - 37b0c38a-...png    (a dependency image)
#   - 6f6f1fc2-...png    (a dependency image)
# and rewrites the Overview / zh-cn / de-de sheets accordingly.

$wb = $excel.ActiveWorkbook

$mdFile   = "19a97789-d66c-41b2-93f9-dea7c2e918c6.md"
$png1     = "37b0c38a-2196-416b-854e-0ad4b7522bac.png"
$png2     = "6f6f1fc2-062f-4c9c-b775-b83e42ff984d.png"
$cfgFile  = ".localization-config"

$readyStatus  = "Ready for handoff"
$notLocalized = "Not to be localized"

$zhXlf   = "19a97789-d66c-41b2-93f9-dea7c2e918c6.72fe77127c097c5be7876c150f8b947c5f62c4da.zh-cn.xlf"
$deXlf   = "19a97789-d66c-41b2-93f9-dea7c2e918c6.72fe77127c097c5be7876c150f8b947c5f62c4da.de-de.xlf"
$png1Dep = "0aec4e24a25c083d3542499701e86068b665678d.png"
$png2Dep = "89426caa34cf234a4a5819f23d3ca2c32c3ad4bd.png"
$depFrom = "e2e\19a97789-d66c-41b2-93f9-dea7c2e918c6.md"

$zhHandoffDt = "2016-03-03 13:18:20"
$deHandoffDt = "2016-03-03 13:18:33"
$epochDt     = "0001-01-01 00:00:00"

$include      = "Include"
$isDependency = "IsDependency"
$ignored      = "Ignored"

$repoBase   = "https://github.com/OpenLocalizationTest/oltest/blob/c162ca4220cfb9c8290489f76192d7792c135d0c"
$zhHtBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d9d5df20bf6ccf9e811473a3ce47023f113236f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"
$deHtBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d7213aeb5ea17aef88c5e81059c3ade7aa8e37d2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"

function Set-RowValues($ws, $row, $values) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

# ---------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Hyperlinks.Delete()

Set-RowValues $ws1 2 @{ "A" = $mdFile;  "B" = $readyStatus;  "C" = $readyStatus }
Set-RowValues $ws1 3 @{ "A" = $png1;    "B" = $readyStatus;  "C" = $readyStatus }
Set-RowValues $ws1 4 @{ "A" = $png2;    "B" = $readyStatus;  "C" = $readyStatus }
Set-RowValues $ws1 5 @{ "A" = $cfgFile; "B" = $notLocalized; "C" = $notLocalized }

$ws1.Hyperlinks.Add($ws1.Range("A2"), "$repoBase/e2e/$mdFile", "", "", $mdFile) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$repoBase/e2e/$png1", "", "", $png1) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$repoBase/e2e/$png2", "", "", $png2) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "$repoBase/$cfgFile", "", "", $cfgFile) | Out-Null

# ---------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Hyperlinks.Delete()

# Clear the columns (E/F) that no longer hold data in rows 2-3
$ws2.Range("E2:F3").ClearContents()

Set-RowValues $ws2 2 @{ "A"=$mdFile; "B"=$readyStatus; "C"=$zhXlf;   "D"=$zhHandoffDt; "G"=$epochDt; "H"=$include }
Set-RowValues $ws2 3 @{ "A"=$png1;   "B"=$readyStatus; "C"=$png1Dep; "D"=$zhHandoffDt; "G"=$epochDt; "H"=$isDependency; "I"=$depFrom }
Set-RowValues $ws2 4 @{ "A"=$png2;   "B"=$readyStatus; "C"=$png2Dep; "D"=$zhHandoffDt; "G"=$epochDt; "H"=$isDependency; "I"=$depFrom }
Set-RowValues $ws2 5 @{ "A"=$cfgFile; "B"=$notLocalized; "D"=$epochDt; "G"=$epochDt; "H"=$ignored }

$ws2.Hyperlinks.Add($ws2.Range("A2"), "$repoBase/e2e/$mdFile", "", "", $mdFile) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$zhHtBase/$zhXlf", "", "", $zhXlf) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$repoBase/e2e/$png1", "", "", $png1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "$zhHtBase/$png1Dep", "", "", $png1Dep) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "$repoBase/e2e/$png2", "", "", $png2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "$zhHtBase/$png2Dep", "", "", $png2Dep) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "$repoBase/$cfgFile", "", "", $cfgFile) | Out-Null

# ---------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Hyperlinks.Delete()

$ws3.Range("E2:F3").ClearContents()

Set-RowValues $ws3 2 @{ "A"=$mdFile; "B"=$readyStatus; "C"=$deXlf;   "D"=$deHandoffDt; "G"=$epochDt; "H"=$include }
Set-RowValues $ws3 3 @{ "A"=$png1;   "B"=$readyStatus; "C"=$png1Dep; "D"=$deHandoffDt; "G"=$epochDt; "H"=$isDependency; "I"=$depFrom }
Set-RowValues $ws3 4 @{ "A"=$png2;   "B"=$readyStatus; "C"=$png2Dep; "D"=$deHandoffDt; "G"=$epochDt; "H"=$isDependency; "I"=$depFrom }
Set-RowValues $ws3 5 @{ "A"=$cfgFile; "B"=$notLocalized; "D"=$epochDt; "G"=$epochDt; "H"=$ignored }

$ws3.Hyperlinks.Add($ws3.Range("A2"), "$repoBase/e2e/$mdFile", "", "", $mdFile) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$deHtBase/$deXlf", "", "", $deXlf) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$repoBase/e2e/$png1", "", "", $png1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "$deHtBase/$png1Dep", "", "", $png1Dep) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "$repoBase/e2e/$png2", "", "", $png2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "$deHtBase/$png2Dep", "", "", $png2Dep) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "$repoBase/$cfgFile", "", "", $cfgFile) | Out-Null
